$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture existing style templates before wiping the sheet ---
# D18 currently carries style index 1 (red "x"/flag style); A2 carries style index 2 (plain heading style).
$ws.Cells.Item(18,4).Copy()
$ws.Cells.Item(200,1).PasteSpecial(-4122)
$ws.Cells.Item(2,1).Copy()
$ws.Cells.Item(200,2).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Wipe existing sheet content + formatting in the old used range ---
$ws.Range("A1:G65").Clear()

# --- Write every cell value for the new layout ---
$ws.Range("A1").Value = 'VNFD'
$ws.Range("A2").Value = 'id'
$ws.Range("A3").Value = 'provider'
$ws.Range("A4").Value = 'product-name'
$ws.Range("A5").Value = 'software-version'
$ws.Range("A6").Value = 'version'
$ws.Range("A7").Value = 'product-info-name'
$ws.Range("A8").Value = 'product-info-description'
$ws.Range("A9").Value = 'vnfm-info'
$ws.Range("A10").Value = 'vdu []'
$ws.Range("B11").Value = 'id'
$ws.Range("B12").Value = 'name'
$ws.Range("B13").Value = 'description'
$ws.Range("B14").Value = 'int-cpd []'
$ws.Range("D15").Value = 'id'
$ws.Range("D16").Value = 'int-virtual-link-desc'
$ws.Range("D17").Value = 'layer-protocol'
$ws.Range("D18").Value = 'additional-sol1-parameters'
$ws.Range("E18").Value = 'x'
$ws.Range("F19").Value = 'allowed-address-variable'
$ws.Range("G19").Value = 'x'
$ws.Range("F20").Value = 'security-group-variable'
$ws.Range("G20").Value = 'x'
$ws.Range("B22").Value = 'boot-order'
$ws.Range("D22").Value = 'key'
$ws.Range("E22").Value = 'x'
$ws.Range("D23").Value = 'value'
$ws.Range("B24").Value = 'virtual-compute-desc'
$ws.Range("C24").Value = 'x'
$ws.Range("B25").Value = 'virtual-storage-desc'
$ws.Range("C25").Value = 'x'
$ws.Range("B26").Value = 'sw-image-desc'
$ws.Range("C26").Value = 'x'
$ws.Range("B27").Value = 'day0 []'
$ws.Range("C27").Value = 'x'
$ws.Range("D27").Value = 'Missing All node of this list'
$ws.Range("B28").Value = 'device-type'
$ws.Range("C28").Value = 'x'
$ws.Range("A30").Value = 'sw-image-desc []'
$ws.Range("B30").Value = 'id'
$ws.Range("D30").Value = 'value need to change'
$ws.Range("B31").Value = 'name'
$ws.Range("D31").Value = 'value need to change'
$ws.Range("B32").Value = 'image-name-variable'
$ws.Range("D32").Value = 'we can always populate this field'
$ws.Range("B33").Value = 'version'
$ws.Range("B34").Value = 'checksum'
$ws.Range("B35").Value = 'container-format'
$ws.Range("B36").Value = 'disk-format'
$ws.Range("B37").Value = 'min-disk'
$ws.Range("B38").Value = 'size'
$ws.Range("B39").Value = 'image'
$ws.Range("A40").Value = 'virtual-compute-descriptor []'
$ws.Range("B41").Value = 'id'
$ws.Range("D41").Value = 'value might need to change'
$ws.Range("B42").Value = 'flavor-name-variable'
$ws.Range("B43").Value = 'virtual-cpu'
$ws.Range("B44").Value = 'virtual-memory'
$ws.Range("A45").Value = 'virtual-storage-descriptor []'
$ws.Range("B46").Value = 'id'
$ws.Range("B47").Value = 'type-of-storage'
$ws.Range("B48").Value = 'size-of-storage'
$ws.Range("B49").Value = 'sw-image-desc'
$ws.Range("A50").Value = 'int-virtual-link-desc[]'
$ws.Range("A51").Value = 'ext-cpd []'
$ws.Range("A52").Value = 'configurable-properties'
$ws.Range("A53").Value = 'df'
$ws.Range("B54").Value = 'id'
$ws.Range("B55").Value = 'description'
$ws.Range("B56").Value = 'vdu-profile []'
$ws.Range("D57").Value = 'id'
$ws.Range("D58").Value = 'min-number-of-instances'
$ws.Range("D59").Value = 'max-number-of-instances'
$ws.Range("D60").Value = 'affinity-or-anti-affinity-group []'
$ws.Range("E60").Value = 'x'
$ws.Range("B61").Value = 'instantiation-level []'
$ws.Range("C61").Value = 'x'
$ws.Range("D61").Value = 'Missing All node of this list'
$ws.Range("B62").Value = 'scaling-aspect []'
$ws.Range("C62").Value = 'x'
$ws.Range("D62").Value = 'Missing All node of this list'
$ws.Range("B63").Value = 'affinity-or-anti-affinity-group []'
$ws.Range("C63").Value = 'x'
$ws.Range("D63").Value = 'Missing All node of this list'

# --- Re-apply style index 1 (red) to the cells that need it ---
$s1Cells = @("D18", "F19", "F20", "D22", "B24", "B25", "B26", "B27", "B28", "D60", "B61", "B62", "B63")
foreach ($a in $s1Cells) {
    $ws.Cells.Item(200,1).Copy()
    $ws.Range($a).PasteSpecial(-4122)
}

# --- Re-apply style index 2 (plain heading) to the cells that need it ---
$s2Cells = @("A2", "A3", "A4", "A5", "A6", "A7", "A8", "A9", "A10", "D16", "A51")
foreach ($a in $s2Cells) {
    $ws.Cells.Item(200,2).Copy()
    $ws.Range($a).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# --- Remove the temporary style-holder cells so they do not widen UsedRange/dimension ---
$ws.Cells.Item(200,1).Clear()
$ws.Cells.Item(200,2).Clear()

# --- Match the final selection shown in the saved file ---
$ws.Range("B61").Select()
